# Add a new sheet "test resistenze" at the end of the workbook with a table
# of expected vs measured resistance values (measured with an ohmmeter).
#
# We copy the existing "VRI" sheet first so the new sheet inherits the same
# sheetFormatPr / pageSetup / pageMargins / printOptions baseline used by the
# rest of this workbook, then we overwrite its contents with the new data.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("VRI")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy([System.Reflection.Missing]::Value, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "test resistenze"

# ---- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "R attesa"
$ws.Range("B1").Value = "exponent"
$ws.Range("C1").Value = "Rmisurata"
$ws.Range("D1").Value = "NOTA: aggiungere 0.3 ohm"
$ws.Range("D2").Value = "(misurati con ohmmetro)"

# ---- Data rows --------------------------------------------------------
# columns: A = R attesa, B = exponent, C = Rmisurata
$rows = @(
    @(2, 1, 1, 1.3),
    @(3, 2, 1, 2.3),
    @(4, 3, 1, 3.5),
    @(5, 4, 1, 4.5),
    @(6, 5, 1, 5.4),
    @(7, 6, 1, 6.3),
    @(8, 7, 1, 7.5),
    @(9, 8, 1, 8.4),
    @(10, 9, 1, 9.4),
    @(11, 10, 1, 10.3),
    @(12, 1, 3, 0.999),
    @(13, 2, 3, 1.996),
    @(14, 3, 3, 2.992),
    @(15, 4, 3, 4.009),
    @(16, 5, 3, 5.006),
    @(17, 6, 3, 6.004),
    @(18, 7, 3, 7.01),
    @(19, 8, 3, 8.01),
    @(20, 9, 3, 9),
    @(21, 10, 3, 10),
    @(22, 100, 3, 100.1),
    @(23, 200, 3, 202.5),
    @(24, 300, 3, 303.5),
    @(25, 400, 3, 408),
    @(26, 500, 3, 508),
    @(27, 600, 3, 610),
    @(28, 0.7, 6, 0.712),
    @(29, 0.8, 6, 0.811),
    @(30, 0.9, 6, 0.913),
    @(31, 1, 6, 1.012),
    @(32, 2, 6, 2.022),
    @(33, 3, 6, 2.991),
    @(34, 4, 6, 4.035),
    @(35, 5, 6, 5.045),
    @(36, 6, 6, 6.055),
    @(37, 7, 6, 7.02),
    @(38, 8, 6, 8.03),
    @(39, 9, 6, 9.04),
    @(40, 10, 6, 10.05)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}

# ---- Number formats ---------------------------------------------------
$ws.Range("C2:C11").NumberFormat = "0.00"
$ws.Range("C12:C21,C28:C30").NumberFormat = "0.000"
$ws.Range("C22:C27").NumberFormat = "0.0"
$ws.Range("A28:A30").NumberFormat = "0.0"

# ---- Column width (column D holds the long notes) ---------------------
$ws.Columns.Item(4).ColumnWidth = 26.02

# ---- Restore page setup (writing to cells resets fit-to-page flags) ---
$ps = $ws.PageSetup
$ps.PrintHeadings = $false
$ps.PrintGridlines = $false
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1

$ws.Range("A1").Select()
